# Commit message: "Fruta / hortaliza, semanal"
#
# A new weekly price-report row is inserted at row 559 of the single
# worksheet, pushing the previously-existing rows 559-683 down to 560-684
# (dimension grows from A1:T683 to A1:T684).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 559; everything below shifts down one.
$ws.Rows(559).Insert()

# Populate the new row with the latest observation.
$ws.Cells.Item(559, 1).Value  = 4
$ws.Cells.Item(559, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(559, 3).Value  = "Los Lagos"
$ws.Cells.Item(559, 4).Value  = 45204
$ws.Cells.Item(559, 5).Value  = 10
$ws.Cells.Item(559, 6).Value  = "Fruta"
$ws.Cells.Item(559, 7).Value  = 100102
$ws.Cells.Item(559, 8).Value  = "Cítricos"
$ws.Cells.Item(559, 9).Value  = 100102006
$ws.Cells.Item(559, 10).Value = "Pomelo"
$ws.Cells.Item(559, 11).Value = "Start Ruby"
$ws.Cells.Item(559, 12).Value = "Primera"
$ws.Cells.Item(559, 13).Value = 100
$ws.Cells.Item(559, 14).Value = 15000
$ws.Cells.Item(559, 15).Value = 15000
$ws.Cells.Item(559, 16).Value = 15000
$ws.Cells.Item(559, 17).Value = "`$/caja 14 kilos empedrada"
$ws.Cells.Item(559, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(559, 19).Value = 1071
$ws.Cells.Item(559, 20).Value = 14
